$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 465, pushing existing rows 465-502 down to 466-503.
$ws.Rows.Item(465).Insert()

# Populate the newly inserted row 465 with a new Ajo price record
# (same shape as neighboring rows; date + volume are the new data points).
$ws.Cells.Item(465, 1).Value = 8
$ws.Cells.Item(465, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(465, 3).Value = "Coquimbo"
$ws.Cells.Item(465, 4).Value = 45132
$ws.Cells.Item(465, 5).Value = 4
$ws.Cells.Item(465, 6).Value = 100112003
$ws.Cells.Item(465, 7).Value = "Ajo"
$ws.Cells.Item(465, 8).Value = "Chino"
$ws.Cells.Item(465, 9).Value = "Primera"
$ws.Cells.Item(465, 10).Value = 300
$ws.Cells.Item(465, 11).Value = 19000
$ws.Cells.Item(465, 12).Value = 20000
$ws.Cells.Item(465, 13).Value = 19500
$ws.Cells.Item(465, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(465, 15).Value = "China"
$ws.Cells.Item(465, 16).Value = 1950
$ws.Cells.Item(465, 17).Value = 10
$ws.Cells.Item(465, 18).Value = "Hortaliza"
